# Auto-generated edit script for ZBP_06_home_office.xlsx
# Adds two new weekly-survey columns (25.-31.10.2021 and 1.-7.11.2021)
# to both worksheets ("data" and "pocetR"), and refreshes the
# "aktualizace" date in the trailing label rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# ---------------------------------------------------------------
# Sheet "data": new columns BP (25.-31. 10. 2021) and BQ (1.-7. 11. 2021)
# ---------------------------------------------------------------

# Header row (row 1) - copy formatting from the previous header cell (BO1)
# so the new header cells keep the same bold/centered/bordered style.
$ws1.Range("BO1").Copy($ws1.Range("BP1"))
$ws1.Range("BP1").Value = "25.–31. 10. 2021"

$ws1.Range("BO1").Copy($ws1.Range("BQ1"))
$ws1.Range("BQ1").Value = "1.–7. 11. 2021"

# Data rows 2-77
$bp1 = @(0.74, 0.09, 0.07000000000000001, 0.1, 0.75, 0.08, 0.07000000000000001, 0.1, 0.49, 0.16, 0.33, 0.02, 0.8, 0.07000000000000001, 0.015, 0.115, 0.84, 0.06, 0.02, 0.08, 0.68, 0.18, 0.04, 0.1, 0.76, 0.09, 0.06, 0.09, 0.79, 0.07000000000000001, 0.03, 0.11, 0.65, 0.1, 0.15, 0.1, 0.64, 0.15, 0.12, 0.09, 0.72, 0.13, 0.1, 0.05, 0.79, 0.07000000000000001, 0.05, 0.09, 0.61, 0.12, 0.1, 0.17, 0.73, 0.07000000000000001, 0.07000000000000001, 0.13, 0.74, 0.09, 0.09, 0.08, 0.74, 0.1, 0.07000000000000001, 0.09, 0.79, 0.07000000000000001, 0.03, 0.11, 0.75, 0.09, 0.08, 0.08, 0.64, 0.14, 0.13, 0.09)
$bq1 = @(0.76, 0.08, 0.07000000000000001, 0.09, 0.76, 0.08, 0.08, 0.08, 0.52, 0.11, 0.33, 0.04, 0.84, 0.05, 0.015, 0.095, 0.8100000000000001, 0.06, 0.015, 0.115, 0.79, 0.1, 0.05, 0.06, 0.79, 0.07000000000000001, 0.06, 0.08, 0.8100000000000001, 0.06, 0.03, 0.1, 0.66, 0.07000000000000001, 0.16, 0.11, 0.67, 0.13, 0.14, 0.06, 0.73, 0.13, 0.1, 0.04, 0.82, 0.04, 0.05, 0.09, 0.63, 0.1, 0.12, 0.15, 0.74, 0.08, 0.07000000000000001, 0.11, 0.78, 0.05, 0.09, 0.08, 0.76, 0.09, 0.07000000000000001, 0.08, 0.79, 0.06, 0.03, 0.12, 0.78, 0.06, 0.08, 0.08, 0.6899999999999999, 0.13, 0.12, 0.06)

for ($i = 0; $i -lt $bp1.Length; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 68).Value = $bp1[$i]
    $ws1.Cells.Item($r, 69).Value = $bq1[$i]
}

# Row 78 label - bump the "aktualizace" date
$ws1.Range("A78").Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 18. 11. 2021"

# ---------------------------------------------------------------
# Sheet "pocetR": new columns BO (25.-31. 10. 2021) and BP (1.-7. 11. 2021)
# ---------------------------------------------------------------

$ws2.Range("BN1").Copy($ws2.Range("BO1"))
$ws2.Range("BO1").Value = "25.–31. 10. 2021"

$ws2.Range("BN1").Copy($ws2.Range("BP1"))
$ws2.Range("BP1").Value = "1.–7. 11. 2021"

# Data rows 2-20
$bo2 = @(955, 258, 90, 254, 134, 83, 471, 232, 118, 134, 237, 540, 178, 267, 236, 452, 353, 373, 229)
$bp2 = @(955, 258, 90, 254, 134, 83, 471, 232, 118, 134, 237, 540, 178, 267, 236, 452, 353, 373, 229)

for ($i = 0; $i -lt $bo2.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 67).Value = $bo2[$i]
    $ws2.Cells.Item($r, 68).Value = $bp2[$i]
}

# Row 21 label - bump the "aktualizace" date, and extend the trailing
# blank placeholder cells to the two new columns (matching columns
# B21:BN21, which are empty text cells). Copying the existing blank
# BN21 cell brings the two new cells into existence without altering
# any style/number-format tables.
$ws2.Range("A21").Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 18. 11. 2021"

$ws2.Range("BN21").Copy($ws2.Range("BO21"))
$ws2.Range("BN21").Copy($ws2.Range("BP21"))

Write-Output "Added columns BP/BQ to 'data' and BO/BP to 'pocetR'."
